# Append two new rows of data (rows 184 and 185) to the active sheet,
# mirroring the result of re-running the R script that produces this
# workbook's data (new trading days appended at the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 184 ---------------------------------------------------------
$r = 184
$ws.Cells.Item($r, 1).Value = 45506.2916666667      # date
$ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Cells.Item($r, 2).Value = 0                     # volume
$ws.Cells.Item($r, 3).Value = 8.25                  # high
$ws.Cells.Item($r, 4).Value = 8.25                  # low
$ws.Cells.Item($r, 5).Value = 8.25                  # open
$ws.Cells.Item($r, 6).Value = 8.25                  # close
$ws.Cells.Item($r, 7).Value = "8.25"                # adj_close (text)
$ws.Cells.Item($r, 8).Value = "VARV.MI"              # ticker

# --- Row 185 ---------------------------------------------------------
$r = 185
$ws.Cells.Item($r, 1).Value = 45509.5947569444      # date
$ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Cells.Item($r, 2).Value = 3400                  # volume
$ws.Cells.Item($r, 3).Value = 8.10000038146973      # high
$ws.Cells.Item($r, 4).Value = 7.90000009536743      # low
$ws.Cells.Item($r, 5).Value = 8.10000038146973      # open
$ws.Cells.Item($r, 6).Value = 7.94999980926514      # close
$ws.Cells.Item($r, 7).Value = "7.94999980926514"    # adj_close (text)
$ws.Cells.Item($r, 8).Value = "VARV.MI"              # ticker
